$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Numeric-looking text values are prefixed with an apostrophe so Excel keeps
# them as text (matching the original inline-string cell contents) instead of
# reinterpreting them as numbers.

$ws.Range("D2").Value = "72.122.79"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.677.50"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'600.06"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "'175.41"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "2.676.86"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "'5.00"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "3.171.91"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "71.963.46"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "'26.31"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "2.677.48"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +5.34%  "
$ws.Range("D20").Value = "'8.21"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").Value = "'373.46"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "'72.08"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'4.36"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Value = "'9.82"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "2.816.65"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "0.0₃0975"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'8.10"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'502.23"
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'161.79"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "'19.09"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'0.111"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").Value = "'0.334"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "'156.46"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").Value = "'3.74"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  -1.61%  "
